# Update countries & provincias Spain
# - Swap four pairs of adjacent country names (their underlying order in the
#   shared-string table was swapped upstream; net effect is the two rows'
#   displayed country names trade places).
# - Refresh the "Datos actualizados" timestamp in A1.
# - Refresh the numeric COVID stat columns (B:H) for the rows whose figures
#   moved between the 16:50 and 18:07 snapshots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Country name swaps (row A-column text) ----
$ws.Range("A41").Value = "Republica Dominicana"
$ws.Range("A42").Value = "Panama"

$ws.Range("A134").Value = "Mozambique"
$ws.Range("A135").Value = "Nueva Zelanda"

$ws.Range("A141").Value = "Liberia"
$ws.Range("A142").Value = "Niger"

$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# ---- Updated "last refreshed" timestamp ----
$ws.Range("A1").Value = "Datos actualizados a 22 de Julio de 2020 a las 18:07"

# ---- Updated statistic columns (Casos totales, Nuevos casos, Casos activos,
#      Recuperados, Casos criticos, Muertes hoy, Muertes) ----

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4046552
$ws.Range("C4").Value = 17983
$ws.Range("D4").Value = 1889285
$ws.Range("E4").Value = 2011996
$ws.Range("G4").Value = 318
$ws.Range("H4").Value = 145271

# Row 6 - India
$ws.Range("B6").Value = 1220433
$ws.Range("C6").Value = 26348
$ws.Range("D6").Value = 772488
$ws.Range("E6").Value = 418414
$ws.Range("G6").Value = 761
$ws.Range("H6").Value = 29531

# Row 11 - Chile
$ws.Range("D11").Value = 309241
$ws.Range("E11").Value = 16720
$ws.Range("G11").Value = 45
$ws.Range("H11").Value = 8722

# Row 13 - Reino Unido
$ws.Range("B13").Value = 296377
$ws.Range("C13").Value = 560
$ws.Range("G13").Value = 79
$ws.Range("H13").Value = 45501

# Row 17 - Italia
$ws.Range("B17").Value = 245032
$ws.Range("C17").Value = 280
$ws.Range("D17").Value = 197628
$ws.Range("E17").Value = 12322
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 35082

# Row 21 - Alemania
$ws.Range("D21").Value = 188600
$ws.Range("E21").Value = 6373

# Row 26 - Irak
$ws.Range("B26").Value = 99865
$ws.Range("C26").Value = 2706
$ws.Range("D26").Value = 67147
$ws.Range("E26").Value = 28676
$ws.Range("G26").Value = 92
$ws.Range("H26").Value = 4042

# Row 30 - Suecia
$ws.Range("B30").Value = 78504
$ws.Range("C30").Value = 132
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = 5667

# Row 41 - Republica Dominicana (after swap)
$ws.Range("B41").Value = 56043
$ws.Range("C41").Value = 1246
$ws.Range("D41").Value = 26466
$ws.Range("E41").Value = 28572
$ws.Range("G41").Value = 6
$ws.Range("H41").Value = 1005

# Row 42 - Panama (after swap)
$ws.Range("B42").Value = 55153
$ws.Range("D42").Value = 30075
$ws.Range("E42").Value = 23919
$ws.Range("H42").Value = 1159

# Row 45 - Portugal
$ws.Range("B45").Value = 49150
$ws.Range("C45").Value = 252
$ws.Range("D45").Value = 33999
$ws.Range("E45").Value = 13449
$ws.Range("G45").Value = 5
$ws.Range("H45").Value = 1702

# Row 71
$ws.Range("B71").Value = 14448
$ws.Range("C71").Value = 124
$ws.Range("D71").Value = 9075
$ws.Range("E71").Value = 5009
$ws.Range("G71").Value = 4
$ws.Range("H71").Value = 364

# Row 101
$ws.Range("B101").Value = 4358
$ws.Range("C101").Value = 68
$ws.Range("D101").Value = 2463
$ws.Range("E101").Value = 1775
$ws.Range("G101").Value = 3
$ws.Range("H101").Value = 120

# Row 132
$ws.Range("B132").Value = 1640
$ws.Range("C132").Value = 11
$ws.Range("D132").Value = 751
$ws.Range("E132").Value = 431
$ws.Range("G132").Value = 2
$ws.Range("H132").Value = 458

# Row 134 - Mozambique (after swap)
$ws.Range("B134").Value = 1557
$ws.Range("C134").Value = 21
$ws.Range("D134").Value = 523
$ws.Range("E134").Value = 1023
$ws.Range("H134").Value = 11

# Row 135 - Nueva Zelanda (after swap)
$ws.Range("B135").Value = 1555
$ws.Range("D135").Value = 1506
$ws.Range("E135").Value = 27
$ws.Range("H135").Value = 22

# Row 137
$ws.Range("B137").Value = 1394
$ws.Range("C137").Value = 5
$ws.Range("D137").Value = 1108

# Row 140
$ws.Range("B140").Value = 1120
$ws.Range("C140").Value = 7
$ws.Range("D140").Value = 1035
$ws.Range("E140").Value = 74

# Row 141 - Liberia (after swap)
$ws.Range("B141").Value = 1114
$ws.Range("C141").Value = 6
$ws.Range("D141").Value = 592
$ws.Range("E141").Value = 452
$ws.Range("H141").Value = 70

# Row 142 - Niger (after swap)
$ws.Range("B142").Value = 1113
$ws.Range("D142").Value = 1018
$ws.Range("E142").Value = 26
$ws.Range("H142").Value = 69
